$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add two new data rows (3 and 4), duplicating the formatting of row 2 ---
$ws.Range("A2:O2").Copy($ws.Range("A3:O3"))
$ws.Range("A2:O2").Copy($ws.Range("A4:O4"))
$excel.CutCopyMode = 0

# Row 3 values (recurring donation, new receipt)
$ws.Range("A3").Value = "AB188009867"
$ws.Range("B3").Value = "2021-01-06 08:21:37"
$ws.Range("C3").Value = 50
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = "Laurence"
$ws.Range("G3").Value = "Sloan"
$ws.Range("H3").Value = "5701 COACH GATE WYNDE Apt 65"
$ws.Range("I3").Value = "Louisville"
$ws.Range("J3").Value = "KY"
$ws.Range("K3").Value = 40207
$ws.Range("L3").Value = "United States"
$ws.Range("M3").Value = "lsloan0033@yahoo.com"
$ws.Range("N3").Value = 5025523990
$ws.Range("O3").Value = "Email45"

# Row 4 values (another recurring donation, new receipt)
$ws.Range("A4").Value = "AB288009867"
$ws.Range("B4").Value = "2021-02-06 08:21:37"
$ws.Range("C4").Value = 50
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = "Laurence"
$ws.Range("G4").Value = "Sloan"
$ws.Range("H4").Value = "5701 COACH GATE WYNDE Apt 65"
$ws.Range("I4").Value = "Louisville"
$ws.Range("J4").Value = "KY"
$ws.Range("K4").Value = 40207
$ws.Range("L4").Value = "United States"
$ws.Range("M4").Value = "lsloan0033@yahoo.com"
$ws.Range("N4").Value = 5025523990
$ws.Range("O4").Value = "Email45"

# --- Hyperlinks for the donor-email cells in the new rows ---
# (Hyperlinks.Add stamps an Office "Hyperlink" style on the cell; restore the
#  plain text style that the rest of the column uses by re-pasting formats
#  from the original, already-hyperlinked cell.)
$ws.Hyperlinks.Add($ws.Range("M3"), "mailto:lsloan0033@yahoo.com")
$ws.Range("M2").Copy()
$ws.Range("M3").PasteSpecial(-4122)

$ws.Hyperlinks.Add($ws.Range("M4"), "mailto:lsloan0033@yahoo.com")
$ws.Range("M2").Copy()
$ws.Range("M4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Column widths: widen the Date column (B) and the Donor Phone column (N) ---
$ws.Columns.Item(2).ColumnWidth = 17.5
$ws.Columns.Item(14).ColumnWidth = 19.166666666666668

# --- Sheet view: scroll / selection state ---
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("G5").Select()
